# Insert a new "Diet" column (G) into Sheet1, shifting the existing
# Body_Weight / Blood_Glucose columns one slot to the right (G->H, H->I).
#
# Diet is "Low Fat" for the FL (Low-Fat) treatment rows and "High Fat" for
# the HF / SC / TV (High-Fat-base) treatment rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing G:H data right by inserting a new column at G.
$ws.Columns("G:G").Insert()

# Header
$ws.Range("G1").Value = "Diet"

# Body (rows 2-77): first block (FL, rows 2-21) is Low Fat; the rest
# (HF/SC/TV, rows 22-77) are all on the High Fat base diet.
$ws.Range("G2:G21").Value = "Low Fat"
$ws.Range("G22:G77").Value = "High Fat"

# Reset the view: scroll back to the top-left and move the selection.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("M68").Select()
